$d = $word.ActiveDocument

$replacements = @(
    @("2026-02-01 Sunday", "2026-02-02 Monday"),
    @("79×93=7347", "39×48=1872"),
    @("49×19=931", "99×66=6534"),
    @("78×29=2262", "44×56=2464"),
    @("40×91=3640", "87×82=7134"),
    @("37×25=925", "28×15=420"),
    @("43×64=2752", "73×96=7008"),
    @("32×76=2432", "41×37=1517"),
    @("47×67=3149", "94×48=4512"),
    @("97×76=7372", "74×49=3626"),
    @("98×80=7840", "29×86=2494"),
    @("39×73=2847", "17×58=986"),
    @("98×18=1764", "92×70=6440"),
    @("76×53=4028", "76×98=7448"),
    @("61×92=5612", "85×90=7650"),
    @("53×27=1431", "32×92=2944"),
    @("21×28=588", "12×99=1188"),
    @("44×78=3432", "13×70=910"),
    @("71×28=1988", "94×96=9024"),
    @("16×32=512", "26×45=1170"),
    @("26×47=1222", "72×58=4176"),
    @("93×59=5487", "97×69=6693"),
    @("19×88=1672", "55×30=1650"),
    @("83×47=3901", "30×25=750"),
    @("47×68=3196", "66×44=2904"),
    @("99×20=1980", "15×59=885")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
